$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap header row (C2:G2) and row labels (B3:B7) between the
#     "Cx" names and the "wrt Cx" names, and reset the pairwise
#     comparison matrix's upper-triangle to all 1s. ---

# Row 2 headers: now use the "wrt Cx" labels
$ws.Range("C2").Value = "wrt C1"
$ws.Range("D2").Value = "wrt C2"
$ws.Range("E2").Value = "wrt C3"
$ws.Range("F2").Value = "wrt C4"
$ws.Range("G2").Value = "wrt C5"

# Column B row labels: now use the plain "Cx" labels
$ws.Range("B3").Value = "C1"
$ws.Range("B4").Value = "C2"
$ws.Range("B5").Value = "C3"
$ws.Range("B6").Value = "C4"
$ws.Range("B7").Value = "C5"

# Reset the comparison values (upper triangle) to 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1

$ws.Range("G5").Value = 1

$ws.Range("G6").Value = 1

# --- Remove the picture / drawing object from the sheet ---
for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete()
}
